# Update metadata sheet "16.b.1.1": indicator code/text refresh, contact info
# refresh (new organization sub-unit, phone, site, e-mail), and tidy up a
# couple of cells whose text had trailing blank lines. Also nudges a few row
# heights / the print scale / margins to match the re-saved layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Информация об индикаторе -------------------------------------------------
# B1 sits next to the "1. Информация об индикаторе" title and now just holds the
# short indicator code instead of the old "16.b.1.1" sub-code.
$ws.Range("B1").Value = "16.b.1"

# B4 = Индикатор: refreshed wording of the 16.b.1 indicator description.
$ws.Range("B4").Value = "16.b.1 Доля лиц, сообщивших о том, что в последние 12 месяцев они лично подвергались дискриминации или преследованиям на основаниях, дискриминация по которым запрещена в соответствии с международными стандартами в области прав человека"

# --- 2. Информация об организации -------------------------------------------------
# B6 = Организация: "Отдел" -> "Управление" статистики домашних хозяйств.
$ws.Range("B6").Value = "Национальный статистический комитет КР (Управление статистики домашних хозяйств)"

# B7 = Контактное лицо (unchanged, left as-is).

# B8 = Электронная почта контактного лица
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com "

# B9 = Телефон контактного лица
$ws.Range("B9").Value = "0(312) 32 46 55"

# B10 = Сайт организации
$ws.Range("B10").Value = "www.stat.gov.kg "

# --- 4. Источники данных и методы сбора -------------------------------------------
# B17 = Методы сбора данных: same wording, just without the trailing blank line.
$ws.Range("B17").Value = "В опросах КОМП используется персональное интервью с использованием компьютеров (Computer-Assisted Personal Interview, CAPI). Приложение для сбора данных, в том числе специальная платформа управления данными КОМП, основано на программном обеспечении CSPro (Census and Survey Processing System) версии 6.3.                                                                  Процедуры и стандартные программы, разработанные в рамках глобальной программы КОМП, были адаптированы к окончательной версии вопросника КОМП Кыргызстан, 2018 года и использовались в опросе. `nДанные собирались с использованием планшетных компьютеров под управлением операционной системы Windows 8 с использованием приложения Bluetooth для операций на местах, что позволяло передавать задания и заполненные вопросники с/на планшеты руководителя и интервьюера."

# --- 5. Метод расчета и другие методологические основы ----------------------------
# B21 = Комментарий и ограничения: same wording, just without the trailing blank line.
$ws.Range("B21").Value = "Руководители команд отвечали за ежедневный мониторинг работ на местах. В одном домохозяйстве из кластера проводился обязательный повторный опрос. Проводились ежедневные наблюдения за навыками и результатами интервьюера. В ходе работ на местах каждую команду несколько раз посещали руководители обследования, а также организовывались визиты на места членов команды КОМП ЮНИСЕФ.`nВ ходе полевых работ каждую неделю создавались проверочные таблицы для анализа и работы с командами на местах, которые представляли собой адаптированные версии стандартных таблиц, созданных Программой КОМП."

# --- row-height touch-ups (content reflow after the edits above) ------------------
$ws.Rows.Item(12).RowHeight = 105.75
$ws.Rows.Item(13).RowHeight = 51.75
$ws.Rows.Item(14).RowHeight = 105
$ws.Rows.Item(17).RowHeight = 189.75
$ws.Rows.Item(21).RowHeight = 146.25
$ws.Rows.Item(23).RowHeight = 62.25

# --- selection state ----------------------------------------------------------------
[void]$ws.Range("B6").Select()

# --- page setup: tighter margins, bigger print scale --------------------------------
# PageSetup margins are in points (1 inch = 72pt); target OOXML inches are
# left/right = 0in, top/bottom = 1.9cm (~0.748in = 53.858pt).
$ws.PageSetup.LeftMargin = 0
$ws.PageSetup.RightMargin = 0
$ws.PageSetup.TopMargin = 53.8582677165354
$ws.PageSetup.BottomMargin = 53.8582677165354
$ws.PageSetup.Zoom = 85
